# Update Price (D) and Volume(1h) (E) columns with the latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.592.32'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '2.449.34'
$ws.Range('E3').Value = '  -2.49%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '548.26'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = '146.87'
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '0.585'
$ws.Range('E8').Value = '  -2.94%  '
$ws.Range('D9').Value = '2.448.38'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('E11').Value = '  -0.06%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '5.40'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '0.351'
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('D14').Value = '26.12'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '2.884.36'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').Value = '0.0000168'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '61.325.77'
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').Value = '2.441.59'
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('D19').Value = '10.92'
$ws.Range('E19').Value = '  -4.24%  '
$ws.Range('D20').Value = '6.95'
$ws.Range('E20').Value = '  -2.66%  '
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -2.65%  '
$ws.Range('D22').Value = '319.35'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D24').Value = '1.89'
$ws.Range('E24').Value = '  +7.02%  '
$ws.Range('D25').Value = '63.66'
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').Value = '0.0₃0985'
$ws.Range('E26').Value = '  -5.46%  '
$ws.Range('D27').Value = '2.561.17'
$ws.Range('E27').Value = '  -3.83%  '
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '537.45'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').Value = '1.48'
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('D31').Value = '8.26'
$ws.Range('E31').Value = '  -4.67%  '
$ws.Range('D32').Value = '7.72'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('E33').Value = '  -4.44%  '
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  -2.39%  '
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('D36').Value = '5.75'
$ws.Range('E36').Value = '  -4.98%  '
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '4.79'
$ws.Range('E38').Value = '  -3.28%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.380'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -1.29%  '
$ws.Range('D40').Value = '18.28'
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = '140.14'
$ws.Range('E42').Value = '  -7.45%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = '40.18'
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('E45').Value = '  -1.68%  '
$ws.Range('D46').Value = '141.72'
$ws.Range('E46').Value = '  -6.15%  '
$ws.Range('D47').Value = '21.98'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = '3.62'
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('D49').Value = '0.0533'
$ws.Range('E49').Value = '  -3.87%  '
$ws.Range('D50').Value = '0.589'
$ws.Range('E50').Value = '  -1.30%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.0930'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -2.51%  '
